$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($cellRef, $text)
    # Values in column D (Price) are stored as plain text in the source
    # workbook (e.g. "30.314.10", "236.35", "0.4710"). Forcing the cell's
    # number format to Text before assignment keeps Excel from
    # re-interpreting numeric-looking strings as actual numbers (which
    # would also destroy meaningful trailing zeros / precision).
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "30.336.14"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.870.63"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "1.002"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
Set-PriceText "D5" "236.51"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
Set-PriceText "D7" "0.4712"
$ws.Range("E7").Value = "  +0.88%  "

# Row 8
Set-PriceText "D8" "0.2897"
$ws.Range("E8").Value = "  +2.30%  "

# Row 9
Set-PriceText "D9" "0.06620"
$ws.Range("E9").Value = "  +1.56%  "

# Row 10
Set-PriceText "D10" "21.68"
$ws.Range("E10").Value = "  +0.11%  "

# Row 11
Set-PriceText "D11" "0.08058"
$ws.Range("E11").Value = "  +1.65%  "

# Row 12
Set-PriceText "D12" "97.27"
$ws.Range("E12").Value = "  -0.05%  "

# Row 13
Set-PriceText "D13" "1.873.33"
$ws.Range("E13").Value = "  +0.48%  "

# Row 14
Set-PriceText "D14" "5.140"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
Set-PriceText "D15" "0.6887"
$ws.Range("E15").Value = "  +1.71%  "

# Row 16
Set-PriceText "D16" "272.00"
$ws.Range("E16").Value = "  -2.34%  "

# Row 17
Set-PriceText "D17" "30.327.07"
$ws.Range("E17").Value = "  +0.14%  "

# Row 18
Set-PriceText "D18" "14.14"
$ws.Range("E18").Value = "  +5.78%  "

# Row 19
Set-PriceText "D19" "0.000007712"
$ws.Range("E19").Value = "  +5.70%  "

# Row 20
Set-PriceText "D20" "1.002"
$ws.Range("E20").Value = "  +0.06%  "

# Row 21
Set-PriceText "D21" "2.118.30"
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
Set-PriceText "D22" "5.315"
$ws.Range("E22").Value = "  -1.22%  "

# Row 23
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
Set-PriceText "D24" "6.218"
$ws.Range("E24").Value = "  +1.09%  "

# Row 25
Set-PriceText "D25" "167.87"
$ws.Range("E25").Value = "  +0.54%  "

# Row 26
Set-PriceText "D26" "9.306"
$ws.Range("E26").Value = "  +1.73%  "

# Row 27
Set-PriceText "D27" "18.97"
$ws.Range("E27").Value = "  -0.51%  "

# Row 28
Set-PriceText "D28" "1.956"
$ws.Range("E28").Value = "  +1.37%  "

# Row 29
Set-PriceText "D29" "1.373"
$ws.Range("E29").Value = "  -0.85%  "

# Row 30
Set-PriceText "D30" "0.09957"
$ws.Range("E30").Value = "  +2.54%  "

# Row 31
Set-PriceText "D31" "4.368"
$ws.Range("E31").Value = "  -0.57%  "

# Row 32
Set-PriceText "D32" "1.465"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33
Set-PriceText "D33" "4.085"
$ws.Range("E33").Value = "  +0.35%  "

# Row 34
Set-PriceText "D34" "0.04703"
$ws.Range("E34").Value = "  -0.71%  "

# Row 35
Set-PriceText "D35" "1.133"
$ws.Range("E35").Value = "  +0.45%  "

# Row 36
Set-PriceText "D36" "0.7024"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37
Set-PriceText "D37" "2.715"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
Set-PriceText "D38" "0.01883"
$ws.Range("E38").Value = "  +1.27%  "

# Row 39
Set-PriceText "D39" "2.648"
$ws.Range("E39").Value = "  +2.82%  "

# Row 40
Set-PriceText "D40" "6.302"
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
Set-PriceText "D41" "73.03"
$ws.Range("E41").Value = "  -1.99%  "

# Row 42
Set-PriceText "D42" "1.962"
$ws.Range("E42").Value = "  +0.26%  "

# Row 43
Set-PriceText "D43" "0.8430"
$ws.Range("E43").Value = "  -0.82%  "

# Row 44
Set-PriceText "D44" "0.4169"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
Set-PriceText "D45" "1.001"
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
Set-PriceText "D46" "103.39"
$ws.Range("E46").Value = "  +0.12%  "

# Row 47 (was Aptos -> now EnergySwap)
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText "D47" "9.291"
$ws.Range("E47").Value = "  -0.36%  "

# Row 48 (was EnergySwap -> now Aptos)
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-PriceText "D48" "7.112"
$ws.Range("E48").Value = "  -0.86%  "

# Row 49
Set-PriceText "D49" "934.12"
$ws.Range("E49").Value = "  -3.45%  "

# Row 50
Set-PriceText "D50" "34.50"

# Row 51
Set-PriceText "D51" "0.05671"
$ws.Range("E51").Value = "  +0.48%  "
